$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Runtime (F) column for rows 343-357 ---
$ws.Cells.Item(343, 6).Value = '15:29:55'
$ws.Cells.Item(344, 6).Value = '15:29:57'
$ws.Cells.Item(345, 6).Value = '15:29:58'
$ws.Cells.Item(346, 6).Value = '15:32:49'
$ws.Cells.Item(347, 6).Value = '15:32:51'
$ws.Cells.Item(348, 6).Value = '15:32:51'
$ws.Cells.Item(349, 6).Value = '15:56:55'
$ws.Cells.Item(350, 6).Value = '15:57:02'
$ws.Cells.Item(351, 6).Value = '15:57:04'
$ws.Cells.Item(352, 6).Value = '17:09:52'
$ws.Cells.Item(353, 6).Value = '17:09:54'
$ws.Cells.Item(354, 6).Value = '17:09:55'
$ws.Cells.Item(355, 6).Value = '17:09:56'
$ws.Cells.Item(356, 6).Value = '15:36:24'
$ws.Cells.Item(357, 6).Value = '15:37:09'

# --- Add Test fail reason (D) column for rows 659-660 ---
$ws.Cells.Item(659, 4).Value = 'Booking completed'
$ws.Cells.Item(660, 4).Value = 'Booking completed'

# --- Touch F358 (empty placeholder cell, matches upstream no-op edit) ---
$ws.Cells.Item(358, 6).Style = "Normal"

# --- Append new rows 661-694 ---
# row 661
$ws.Cells.Item(661, 1).Value = 'Booking_01'
$ws.Cells.Item(661, 2).Value = 'ZA12580363'
$ws.Cells.Item(661, 3).Value = 'Passed'
$ws.Cells.Item(661, 4).Value = 'Booking completed'
$ws.Cells.Item(661, 5).Value = '62d083e5-c481-4995-b365-693246a1aee2'
# row 662
$ws.Cells.Item(662, 1).Value = 'Booking_02'
$ws.Cells.Item(662, 2).Value = 'ZA12580376'
$ws.Cells.Item(662, 3).Value = 'Passed'
$ws.Cells.Item(662, 4).Value = 'Skipped this test case as this test case is not approved to run'
$ws.Cells.Item(662, 5).Value = '5f8224a6-298e-4020-a14c-02414e28e247'
# row 663
$ws.Cells.Item(663, 1).Value = 'Booking_03'
$ws.Cells.Item(663, 2).Value = 'ZA12580383'
$ws.Cells.Item(663, 3).Value = 'Passed'
$ws.Cells.Item(663, 4).Value = 'Skipped this test case as this test case is not approved to run'
$ws.Cells.Item(663, 5).Value = '90c84593-b6b0-473b-987e-cadc69c967fc'
# row 664
$ws.Cells.Item(664, 1).Value = 'Booking_01'
$ws.Cells.Item(664, 2).Value = '-'
$ws.Cells.Item(664, 3).Value = 'Skipped'
$ws.Cells.Item(664, 4).Value = 'Skipped this test case as this test case is not approved to run'
$ws.Cells.Item(664, 5).Value = '-'
# row 665
$ws.Cells.Item(665, 1).Value = 'Booking_02'
$ws.Cells.Item(665, 2).Value = '-'
$ws.Cells.Item(665, 3).Value = 'Skipped'
$ws.Cells.Item(665, 4).Value = 'Skipped this test case as this test case is not approved to run'
$ws.Cells.Item(665, 5).Value = '-'
# row 666
$ws.Cells.Item(666, 1).Value = 'Booking_03'
$ws.Cells.Item(666, 2).Value = '-'
$ws.Cells.Item(666, 3).Value = 'Skipped'
$ws.Cells.Item(666, 4).Value = 'Skipped this test case as this test case is not approved to run'
$ws.Cells.Item(666, 5).Value = '-'
# row 667
$ws.Cells.Item(667, 1).Value = 'Booking_01'
$ws.Cells.Item(667, 2).Value = '-'
$ws.Cells.Item(667, 3).Value = 'Skipped'
$ws.Cells.Item(667, 4).Value = 'Skipped this test case as this test case is not approved to run'
$ws.Cells.Item(667, 5).Value = '-'
# row 668
$ws.Cells.Item(668, 1).Value = 'Booking_02'
$ws.Cells.Item(668, 2).Value = '-'
$ws.Cells.Item(668, 3).Value = 'Skipped'
$ws.Cells.Item(668, 4).Value = 'Skipped this test case as this test case is not approved to run'
$ws.Cells.Item(668, 5).Value = '-'
# row 669
$ws.Cells.Item(669, 1).Value = 'Booking_03'
$ws.Cells.Item(669, 2).Value = '-'
$ws.Cells.Item(669, 3).Value = 'Skipped'
$ws.Cells.Item(669, 4).Value = 'Skipped this test case as this test case is not approved to run'
$ws.Cells.Item(669, 5).Value = '-'
# row 670
$ws.Cells.Item(670, 1).Value = 'Booking_01'
$ws.Cells.Item(670, 2).Value = '-'
$ws.Cells.Item(670, 3).Value = 'Skipped'
$ws.Cells.Item(670, 4).Value = 'Skipped this test case as this test case is not approved to run'
$ws.Cells.Item(670, 5).Value = '-'
# row 671
$ws.Cells.Item(671, 1).Value = 'Booking_02'
$ws.Cells.Item(671, 2).Value = '-'
$ws.Cells.Item(671, 3).Value = 'Skipped'
$ws.Cells.Item(671, 4).Value = 'Booking completed'
$ws.Cells.Item(671, 5).Value = '-'
# row 672
$ws.Cells.Item(672, 1).Value = 'Booking_03'
$ws.Cells.Item(672, 2).Value = '-'
$ws.Cells.Item(672, 3).Value = 'Skipped'
$ws.Cells.Item(672, 4).Value = 'Booking completed'
$ws.Cells.Item(672, 5).Value = '-'
# row 673
$ws.Cells.Item(673, 1).Value = 'Booking_01'
$ws.Cells.Item(673, 2).Value = 'ZA12590999'
$ws.Cells.Item(673, 3).Value = 'Passed'
$ws.Cells.Item(673, 4).Value = 'Booking completed'
$ws.Cells.Item(673, 5).Value = '8e93d771-f692-4bce-bda8-b82ad876378e'
# row 674
$ws.Cells.Item(674, 1).Value = 'Booking_01'
$ws.Cells.Item(674, 2).Value = 'ZA12591021'
$ws.Cells.Item(674, 3).Value = 'Passed'
$ws.Cells.Item(674, 4).Value = 'Booking completed'
$ws.Cells.Item(674, 5).Value = 'f4dc73c5-5626-453f-8580-4d85f76c8ef5'
# row 675
$ws.Cells.Item(675, 1).Value = 'Booking_01'
$ws.Cells.Item(675, 2).Value = 'ZA12591036'
$ws.Cells.Item(675, 3).Value = 'Passed'
$ws.Cells.Item(675, 4).Value = 'Booking completed'
$ws.Cells.Item(675, 5).Value = '66c11f4e-de89-4ad6-957f-cfbcd31f93bf'
# row 676
$ws.Cells.Item(676, 1).Value = 'Booking_02'
$ws.Cells.Item(676, 2).Value = 'ZA12591044'
$ws.Cells.Item(676, 3).Value = 'Passed'
$ws.Cells.Item(676, 4).Value = 'Booking completed'
$ws.Cells.Item(676, 5).Value = '7f1600a1-3125-4af2-a06f-16a2df6c6af0'
# row 677
$ws.Cells.Item(677, 1).Value = 'Booking_03'
$ws.Cells.Item(677, 2).Value = 'ZA12591049'
$ws.Cells.Item(677, 3).Value = 'Passed'
$ws.Cells.Item(677, 4).Value = 'Skipped this test case as this test case is not approved to run'
$ws.Cells.Item(677, 5).Value = '33d89878-1fd1-40c4-bec7-1325a2333c52'
# row 678
$ws.Cells.Item(678, 1).Value = 'Booking_04'
$ws.Cells.Item(678, 2).Value = 'ZA12591058'
$ws.Cells.Item(678, 3).Value = 'Passed'
$ws.Cells.Item(678, 4).Value = 'Skipped this test case as this test case is not approved to run'
$ws.Cells.Item(678, 5).Value = 'fb0f0bb4-1bc5-44f8-9a9a-9bb286438d11'
# row 679
$ws.Cells.Item(679, 1).Value = 'Booking_01'
$ws.Cells.Item(679, 2).Value = '-'
$ws.Cells.Item(679, 3).Value = 'Skipped'
$ws.Cells.Item(679, 4).Value = 'Skipped this test case as this test case is not approved to run'
$ws.Cells.Item(679, 5).Value = '-'
# row 680
$ws.Cells.Item(680, 1).Value = 'Booking_02'
$ws.Cells.Item(680, 2).Value = '-'
$ws.Cells.Item(680, 3).Value = 'Skipped'
$ws.Cells.Item(680, 4).Value = 'Skipped this test case as this test case is not approved to run'
$ws.Cells.Item(680, 5).Value = '-'
# row 681
$ws.Cells.Item(681, 1).Value = 'Booking_03'
$ws.Cells.Item(681, 2).Value = '-'
$ws.Cells.Item(681, 3).Value = 'Skipped'
$ws.Cells.Item(681, 4).Value = 'Booking completed'
$ws.Cells.Item(681, 5).Value = '-'
# row 682
$ws.Cells.Item(682, 1).Value = 'Booking_04'
$ws.Cells.Item(682, 2).Value = '-'
$ws.Cells.Item(682, 3).Value = 'Skipped'
$ws.Cells.Item(682, 4).Value = 'Booking completed'
$ws.Cells.Item(682, 5).Value = '-'
# row 683
$ws.Cells.Item(683, 1).Value = 'Booking_01'
$ws.Cells.Item(683, 2).Value = 'NG12591099'
$ws.Cells.Item(683, 3).Value = 'Passed'
$ws.Cells.Item(683, 4).Value = 'Booking completed'
$ws.Cells.Item(683, 5).Value = '9e3d49c0-1dfa-4c6e-b347-dac11fcb7528'
# row 684
$ws.Cells.Item(684, 1).Value = 'Booking_02'
$ws.Cells.Item(684, 2).Value = 'NG12591109'
$ws.Cells.Item(684, 3).Value = 'Passed'
$ws.Cells.Item(684, 4).Value = 'Booking completed'
$ws.Cells.Item(684, 5).Value = '8bb7403e-740e-4f74-a11c-b9c4453f0389'
# row 685
$ws.Cells.Item(685, 1).Value = 'Booking_04'
$ws.Cells.Item(685, 2).Value = 'NG12591116'
$ws.Cells.Item(685, 3).Value = 'Passed'
$ws.Cells.Item(685, 4).Value = 'Booking completed'
$ws.Cells.Item(685, 5).Value = '86c5ecdb-993a-4e26-8bfa-b4410d3c2fb2'
# row 686
$ws.Cells.Item(686, 1).Value = 'Booking_01'
$ws.Cells.Item(686, 2).Value = 'ZA00109536'
$ws.Cells.Item(686, 3).Value = 'Passed'
$ws.Cells.Item(686, 4).Value = 'Booking completed'
$ws.Cells.Item(686, 5).Value = '3bfa3cbc-d542-4785-ad6f-3b76da7be95b'
# row 687
$ws.Cells.Item(687, 1).Value = 'Booking_02'
$ws.Cells.Item(687, 2).Value = 'ZA00109537'
$ws.Cells.Item(687, 3).Value = 'Passed'
$ws.Cells.Item(687, 4).Value = 'Booking completed'
$ws.Cells.Item(687, 5).Value = '310e573a-c70a-4c91-9035-d0e3435d24b7'
# row 688
$ws.Cells.Item(688, 1).Value = 'Booking_03'
$ws.Cells.Item(688, 2).Value = 'ZA00109538'
$ws.Cells.Item(688, 3).Value = 'Passed'
$ws.Cells.Item(688, 4).Value = 'Skipped this test case as this test case is not approved to run'
$ws.Cells.Item(688, 5).Value = '7a5af4b3-b367-42e2-bdd8-4fe80a88fe92'
# row 689
$ws.Cells.Item(689, 1).Value = 'Booking_04'
$ws.Cells.Item(689, 2).Value = 'ZA00109539'
$ws.Cells.Item(689, 3).Value = 'Passed'
$ws.Cells.Item(689, 4).Value = 'Skipped this test because desired airline: QR was not avaible in result'
$ws.Cells.Item(689, 5).Value = '39c60eba-2979-4e27-a883-f0eefffff3a0'
# row 690
$ws.Cells.Item(690, 1).Value = 'Booking_05'
$ws.Cells.Item(690, 2).Value = '-'
$ws.Cells.Item(690, 3).Value = 'Skipped'
$ws.Cells.Item(690, 4).Value = 'Skipped this test case as this test case is not approved to run'
$ws.Cells.Item(690, 5).Value = '-'
# row 691
$ws.Cells.Item(691, 1).Value = 'Booking_06'
$ws.Cells.Item(691, 2).Value = '-'
$ws.Cells.Item(691, 3).Value = 'Skipped'
$ws.Cells.Item(691, 4).Value = 'Skipped this test because desired airline: QR was not avaible in result'
$ws.Cells.Item(691, 5).Value = '7dbc8555-20a8-40ec-9a05-9a95479a665f'
# row 692
$ws.Cells.Item(692, 1).Value = 'Booking_07'
$ws.Cells.Item(692, 2).Value = '-'
$ws.Cells.Item(692, 3).Value = 'Skipped'
$ws.Cells.Item(692, 5).Value = '-'
# row 693
$ws.Cells.Item(693, 1).Value = 'Booking_08'
$ws.Cells.Item(693, 2).Value = '-'
$ws.Cells.Item(693, 3).Value = 'Skipped'
$ws.Cells.Item(693, 5).Value = '1dab8b91-acd2-498d-a17c-054839bc1afd'
# row 694
$ws.Cells.Item(694, 1).Value = 'Booking_01'
$ws.Cells.Item(694, 2).Value = 'ZA00109624'
